$d = $word.ActiveDocument

$pairs = @(
    @("27÷3=9, 0", "45÷8=5, 5"),
    @("58÷2=29, 0", "33÷9=3, 6"),
    @("22÷3=7, 1", "92÷4=23, 0"),
    @("50÷6=8, 2", "60÷2=30, 0"),
    @("16÷4=4, 0", "14÷8=1, 6"),
    @("74÷9=8, 2", "50÷2=25, 0"),
    @("24÷6=4, 0", "88÷9=9, 7"),
    @("66÷7=9, 3", "99÷4=24, 3"),
    @("19÷9=2, 1", "34÷8=4, 2"),
    @("39÷3=13, 0", "76÷4=19, 0"),
    @("24÷3=8, 0", "84÷9=9, 3"),
    @("55÷5=11, 0", "40÷7=5, 5"),
    @("77÷9=8, 5", "81÷2=40, 1"),
    @("48÷5=9, 3", "26÷9=2, 8"),
    @("99÷2=49, 1", "11÷5=2, 1"),
    @("24÷8=3, 0", "66÷7=9, 3"),
    @("26÷7=3, 5", "68÷2=34, 0"),
    @("96÷4=24, 0", "60÷8=7, 4"),
    @("62÷4=15, 2", "76÷9=8, 4"),
    @("95÷8=11, 7", "90÷2=45, 0"),
    @("29÷4=7, 1", "17÷5=3, 2"),
    @("68÷8=8, 4", "92÷5=18, 2"),
    @("85÷5=17, 0", "13÷9=1, 4"),
    @("11÷4=2, 3", "93÷6=15, 3"),
    @("79÷4=19, 3", "32÷3=10, 2")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
